# Remove the trailing "Ver no Jupiter..." / "© 2020..." footer paragraphs
# (and the blank paragraph that precedes them) that follow the
# "LOB1004: Cálculo II (Requisito fraco)" requirement line.

$d = $word.ActiveDocument

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
$anchor  = "LOB1004: Cálculo II (Requisito fraco)"

# Walk paragraphs from the end towards the start so earlier deletions never
# invalidate the indices of paragraphs still queued for removal.
$paras = @($d.Paragraphs)
$toDelete = New-Object System.Collections.ArrayList

for ($i = 0; $i -lt $paras.Count; $i++) {
    $text = $paras[$i].Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $target1 -or $text -eq $target2) {
        [void]$toDelete.Add($i)
    }
}

# Also queue the blank paragraph immediately following the anchor line.
for ($i = 0; $i -lt $paras.Count; $i++) {
    $text = $paras[$i].Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $anchor) {
        if ($i + 1 -lt $paras.Count) {
            $nextText = $paras[$i + 1].Range.Text.TrimEnd([char]13, [char]7)
            if ($nextText -eq "") {
                [void]$toDelete.Add($i + 1)
            }
        }
        break
    }
}

$sorted = $toDelete | Sort-Object -Descending -Unique
foreach ($idx in $sorted) {
    $paras[$idx].Range.Delete()
}
